$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.715624570846558
$ws.Range("B1").Value = 3.756839752197266
$ws.Range("C1").Value = 2.080960750579834
$ws.Range("D1").Value = 1.454201817512512
$ws.Range("E1").Value = 1.241882562637329
